# Lab Notebook and size images
# - Fill in the previously-missing GMS/Surface sample count (C4 = 7)
# - Move the active selection from D4 to D6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 7

$ws.Range("D6").Select()
